$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit reshuffles the data rows (2,3,4,5,7,8) of the weekly price table.
# Row 6 is left untouched. For each destination row we write the full set of
# columns D, L, M, N, O, P, Q, R, S, T that changed in the diff.

# Row 2 <- old Row 5
$ws.Range("D2").Value = 44159
$ws.Range("L2").Value = "Segunda"
$ws.Range("M2").Value = 200
$ws.Range("N2").Value = 6500
$ws.Range("O2").Value = 7000
$ws.Range("P2").Value = 6750
$ws.Range("Q2").Value = "$/bandeja 12 canastillos 125 gramos"
$ws.Range("R2").Value = "Provincia de Curicó"
$ws.Range("S2").Value = 4500
$ws.Range("T2").Value = 1.5

# Row 3 <- old Row 8
$ws.Range("D3").Value = 44516
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 80
$ws.Range("N3").Value = 3700
$ws.Range("O3").Value = 3800
$ws.Range("P3").Value = 3750
$ws.Range("Q3").Value = "$/kilo"
$ws.Range("R3").Value = "Región del Maule"
$ws.Range("S3").Value = 3750
$ws.Range("T3").Value = 1

# Row 4 <- old Row 3
$ws.Range("D4").Value = 44162
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 100
$ws.Range("N4").Value = 7000
$ws.Range("O4").Value = 7000
$ws.Range("P4").Value = 7000
$ws.Range("Q4").Value = "$/bandeja 12 canastillos 125 gramos"
$ws.Range("R4").Value = "Provincia de Curicó"
$ws.Range("S4").Value = 4667
$ws.Range("T4").Value = 1.5

# Row 5 <- old Row 4
$ws.Range("D5").Value = 44162
$ws.Range("L5").Value = "Segunda"
$ws.Range("M5").Value = 100
$ws.Range("N5").Value = 6500
$ws.Range("O5").Value = 6500
$ws.Range("P5").Value = 6500
$ws.Range("Q5").Value = "$/bandeja 12 canastillos 125 gramos"
$ws.Range("R5").Value = "Provincia de Curicó"
$ws.Range("S5").Value = 4333
$ws.Range("T5").Value = 1.5

# Row 7 <- old Row 2
$ws.Range("D7").Value = 44176
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 300
$ws.Range("N7").Value = 5000
$ws.Range("O7").Value = 6000
$ws.Range("P7").Value = 5500
$ws.Range("Q7").Value = "$/bandeja 12 canastillos 125 gramos"
$ws.Range("R7").Value = "Provincia de Curicó"
$ws.Range("S7").Value = 3667
$ws.Range("T7").Value = 1.5

# Row 8 <- old Row 7
$ws.Range("D8").Value = 44166
$ws.Range("L8").Value = "Primera"
$ws.Range("M8").Value = 200
$ws.Range("N8").Value = 6000
$ws.Range("O8").Value = 6500
$ws.Range("P8").Value = 6250
$ws.Range("Q8").Value = "$/bandeja 12 canastillos 125 gramos"
$ws.Range("R8").Value = "Provincia de Curicó"
$ws.Range("S8").Value = 4167
$ws.Range("T8").Value = 1.5

Write-Host "Applied weekly reshuffle edits"
